$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-09-03 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-09-04 Monday", 2) | Out-Null
$d.Content.Find.Execute("25÷5=5, 0", $true, $true, $false, $false, $false, $true, 1, $false, "66÷9=7, 3", 2) | Out-Null
$d.Content.Find.Execute("87÷3=29, 0", $true, $true, $false, $false, $false, $true, 1, $false, "28÷4=7, 0", 2) | Out-Null
$d.Content.Find.Execute("79÷2=39, 1", $true, $true, $false, $false, $false, $true, 1, $false, "53÷6=8, 5", 2) | Out-Null
$d.Content.Find.Execute("93÷4=23, 1", $true, $true, $false, $false, $false, $true, 1, $false, "80÷3=26, 2", 2) | Out-Null
$d.Content.Find.Execute("34÷6=5, 4", $true, $true, $false, $false, $false, $true, 1, $false, "56÷2=28, 0", 2) | Out-Null
$d.Content.Find.Execute("47÷9=5, 2", $true, $true, $false, $false, $false, $true, 1, $false, "19÷5=3, 4", 2) | Out-Null
$d.Content.Find.Execute("14÷2=7, 0", $true, $true, $false, $false, $false, $true, 1, $false, "62÷9=6, 8", 2) | Out-Null
$d.Content.Find.Execute("43÷8=5, 3", $true, $true, $false, $false, $false, $true, 1, $false, "30÷3=10, 0", 2) | Out-Null
$d.Content.Find.Execute("57÷6=9, 3", $true, $true, $false, $false, $false, $true, 1, $false, "82÷8=10, 2", 2) | Out-Null
$d.Content.Find.Execute("65÷5=13, 0", $true, $true, $false, $false, $false, $true, 1, $false, "78÷4=19, 2", 2) | Out-Null
$d.Content.Find.Execute("60÷2=30, 0", $true, $true, $false, $false, $false, $true, 1, $false, "97÷5=19, 2", 2) | Out-Null
$d.Content.Find.Execute("58÷5=11, 3", $true, $true, $false, $false, $false, $true, 1, $false, "35÷4=8, 3", 2) | Out-Null
$d.Content.Find.Execute("58÷3=19, 1", $true, $true, $false, $false, $false, $true, 1, $false, "80÷7=11, 3", 2) | Out-Null
$d.Content.Find.Execute("83÷7=11, 6", $true, $true, $false, $false, $false, $true, 1, $false, "81÷8=10, 1", 2) | Out-Null
$d.Content.Find.Execute("36÷3=12, 0", $true, $true, $false, $false, $false, $true, 1, $false, "77÷5=15, 2", 2) | Out-Null
$d.Content.Find.Execute("68÷7=9, 5", $true, $true, $false, $false, $false, $true, 1, $false, "56÷7=8, 0", 2) | Out-Null
$d.Content.Find.Execute("59÷2=29, 1", $true, $true, $false, $false, $false, $true, 1, $false, "73÷5=14, 3", 2) | Out-Null
$d.Content.Find.Execute("55÷9=6, 1", $true, $true, $false, $false, $false, $true, 1, $false, "28÷2=14, 0", 2) | Out-Null
$d.Content.Find.Execute("88÷4=22, 0", $true, $true, $false, $false, $false, $true, 1, $false, "32÷3=10, 2", 2) | Out-Null
$d.Content.Find.Execute("56÷9=6, 2", $true, $true, $false, $false, $false, $true, 1, $false, "13÷4=3, 1", 2) | Out-Null
$d.Content.Find.Execute("81÷6=13, 3", $true, $true, $false, $false, $false, $true, 1, $false, "70÷6=11, 4", 2) | Out-Null
$d.Content.Find.Execute("99÷7=14, 1", $true, $true, $false, $false, $false, $true, 1, $false, "16÷7=2, 2", 2) | Out-Null
$d.Content.Find.Execute("14÷9=1, 5", $true, $true, $false, $false, $false, $true, 1, $false, "34÷4=8, 2", 2) | Out-Null
$d.Content.Find.Execute("10÷2=5, 0", $true, $true, $false, $false, $false, $true, 1, $false, "48÷2=24, 0", 2) | Out-Null
$d.Content.Find.Execute("87÷2=43, 1", $true, $true, $false, $false, $false, $true, 1, $false, "43÷9=4, 7", 2) | Out-Null
